$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 (Login to Appium for iPhone) - procedure text in column E:
# "1. Launch Jabber Client-A" -> "1. Launch Appium Client-A"
$cellE8 = $ws.Range("E8")
$cellE8.Value = $cellE8.Value2.Replace("1. Launch Jabber Client-A", "1. Launch Appium Client-A")

# Row 9 (Make a call to Contact) - procedure text in column E:
# "...Appium  Client-A & Jabber Client-B" -> "...Appium  Client-A & Appium Client-B"
$cellE9 = $ws.Range("E9")
$cellE9.Value = $cellE9.Value2.Replace("Appium  Client-A & Jabber Client-B", "Appium  Client-A & Appium Client-B")

# Row 10 (Verify that user is able to make call to contact) - procedure text in column E:
# "1.. Launch Jabber Client-A & " -> "1.. Launch Appuim Client-A & "
$cellE10 = $ws.Range("E10")
$cellE10.Value = $cellE10.Value2.Replace("1.. Launch Jabber Client-A & ", "1.. Launch Appuim Client-A & ")

# Update the active selection/view as it was left after the edit: cell H4 selected,
# scrolled back to the top (no frozen/offset top-left cell).
$ws.Range("H4").Select()
